$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 770220
$ws.Range("I6").Value = 834330
$ws.Range("J6").Value = 900
$ws.Range("K6").Value = 2502990
$ws.Range("L6").Value = 2700
$ws.Range("M6").Value = -2502878
$ws.Range("N6").Value = -2924

$ws.Range("H17").Value = 844148.0600000001
$ws.Range("J17").Value = 2059923.8
$ws.Range("L17").Value = 6179771.4
$ws.Range("N17").Value = -6180107.4

$ws.Range("H32").Value = 3197.375
$ws.Range("I32").Value = 2359
$ws.Range("K32").Value = 2359
$ws.Range("M32").Value = -2033

$ws.Range("H41").Value = 811.7857
$ws.Range("I41").Value = 751.7778
$ws.Range("J41").Value = 919.8
$ws.Range("K41").Value = 751.7778
$ws.Range("L41").Value = 919.8
$ws.Range("M41").Value = -311.7778
$ws.Range("N41").Value = -1799.8

$ws.Range("H62").Value = 2629.3333
$ws.Range("I62").Value = 2629.3333
$ws.Range("K62").Value = 2629.3333
$ws.Range("M62").Value = -2005.3333

$ws.Range("H65").Value = 2629.3333
$ws.Range("I65").Value = 2629.3333
$ws.Range("K65").Value = 13146.6665
$ws.Range("M65").Value = -10026.6665

$ws.Range("H137").Value = 2311
$ws.Range("I137").Value = 3251
$ws.Range("J137").Value = 1997.6666
$ws.Range("K137").Value = 9753
$ws.Range("L137").Value = 5992.9998
$ws.Range("M137").Value = -7203
$ws.Range("N137").Value = -11092.9998

$ws.Range("H138").Value = 7355860
$ws.Range("I138").Value = 1339
$ws.Range("J138").Value = 10207613
$ws.Range("K138").Value = 4017
$ws.Range("L138").Value = 30622839
$ws.Range("M138").Value = 1123
$ws.Range("N138").Value = -30633119

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1808.091
$ws.Range("I2").Value = 1808.091
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = 1808.091
$ws.Range("L2").Value = 0
$ws.Range("M2").Value = -1695.091
$ws.Range("N2").Value = $null

$ws.Range("H32").Value = 5462.0654
$ws.Range("I32").Value = 4983.091
$ws.Range("K32").Value = 4983.091
$ws.Range("M32").Value = -4696.091

$ws.Range("H36").Value = 4013
$ws.Range("I36").Value = 4013
$ws.Range("J36").Value = 0
$ws.Range("K36").Value = 4013
$ws.Range("L36").Value = 0
$ws.Range("M36").Value = -3667
$ws.Range("N36").Value = $null

$ws.Range("H61").Value = 4150.0435
$ws.Range("I61").Value = 3091.5454
$ws.Range("J61").Value = 5120.3335
$ws.Range("K61").Value = 3091.5454
$ws.Range("L61").Value = 5120.3335
$ws.Range("M61").Value = -2879.5454
$ws.Range("N61").Value = -5544.3335

$ws.Range("H116").Value = 1808.091
$ws.Range("I116").Value = 1808.091
$ws.Range("J116").Value = 0
$ws.Range("K116").Value = 1808.091
$ws.Range("L116").Value = 0
$ws.Range("M116").Value = 485.9090000000001
$ws.Range("N116").Value = $null

$ws.Range("H122").Value = 2849.6155
$ws.Range("I122").Value = 2784.4443
$ws.Range("J122").Value = 2996.25
$ws.Range("K122").Value = 8353.332900000001
$ws.Range("L122").Value = 8988.75
$ws.Range("M122").Value = -5903.332900000001
$ws.Range("N122").Value = -13888.75

$ws.Range("H136").Value = 4150.0435
$ws.Range("I136").Value = 3091.5454
$ws.Range("J136").Value = 5120.3335
$ws.Range("K136").Value = 9274.636200000001
$ws.Range("L136").Value = 15361.0005
$ws.Range("M136").Value = -6724.636200000001
$ws.Range("N136").Value = -20461.0005

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1808.091
$ws.Range("I3").Value = 1808.091
$ws.Range("J3").Value = 0
$ws.Range("K3").Value = 1808.091
$ws.Range("L3").Value = 0
$ws.Range("M3").Value = -1694.091
$ws.Range("N3").Value = $null

$ws.Range("H20").Value = 2936.9792
$ws.Range("I20").Value = 2375.6775
$ws.Range("J20").Value = 3960.5293
$ws.Range("K20").Value = 2375.6775
$ws.Range("L20").Value = 3960.5293
$ws.Range("M20").Value = -2128.6775
$ws.Range("N20").Value = -4454.5293

$ws.Range("H134").Value = 1806.6154
$ws.Range("I134").Value = 1748.871
$ws.Range("K134").Value = 5246.613
$ws.Range("M134").Value = -2711.613

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H20").Value = 0
$ws.Range("J20").Value = 0
$ws.Range("L20").Value = 0
$ws.Range("N20").Value = $null

$ws.Range("H30").Value = 0
$ws.Range("J30").Value = 0
$ws.Range("L30").Value = 0
$ws.Range("N30").Value = $null

$ws.Range("H99").Value = 5061.0586
$ws.Range("I99").Value = 5012.5835
$ws.Range("J99").Value = 5177.4
$ws.Range("K99").Value = 5012.5835
$ws.Range("L99").Value = 5177.4
$ws.Range("M99").Value = -3514.5835
$ws.Range("N99").Value = -8173.4

$ws.Range("H126").Value = 5061.0586
$ws.Range("I126").Value = 5012.5835
$ws.Range("J126").Value = 5177.4
$ws.Range("K126").Value = 15037.7505
$ws.Range("L126").Value = 15532.2
$ws.Range("M126").Value = -12567.7505
$ws.Range("N126").Value = -20472.2

$ws.Range("H127").Value = 100722
$ws.Range("J127").Value = 100722
$ws.Range("L127").Value = 100722
$ws.Range("N127").Value = -110642

$ws.Range("H128").Value = 0
$ws.Range("J128").Value = 0
$ws.Range("L128").Value = 0
$ws.Range("N128").Value = $null

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 60536040
$ws.Range("I4").Value = 68607490
$ws.Range("K4").Value = 205822470
$ws.Range("M4").Value = -205822358

$ws.Range("H129").Value = 1727.7778
$ws.Range("I129").Value = 1258.3334
$ws.Range("J129").Value = 2666.6667
$ws.Range("K129").Value = 3775.0002
$ws.Range("L129").Value = 8000.000100000001
$ws.Range("M129").Value = 1224.9998
$ws.Range("N129").Value = -18000.0001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 2478
$ws.Range("I80").Value = 2336.25
$ws.Range("J80").Value = 2640
$ws.Range("K80").Value = 2336.25
$ws.Range("L80").Value = 2640
$ws.Range("M80").Value = -1338.25
$ws.Range("N80").Value = -4636

$ws.Range("H83").Value = 2478
$ws.Range("I83").Value = 2336.25
$ws.Range("J83").Value = 2640
$ws.Range("K83").Value = 11681.25
$ws.Range("L83").Value = 13200
$ws.Range("M83").Value = -6689.25
$ws.Range("N83").Value = -23184

$ws.Range("H126").Value = 18427.867
$ws.Range("J126").Value = 3533.3333
$ws.Range("L126").Value = 10599.9999
$ws.Range("N126").Value = -15539.9999

$ws.Range("H132").Value = 2754.5356
$ws.Range("I132").Value = 2804.923
$ws.Range("K132").Value = 8414.769
$ws.Range("M132").Value = -5884.769

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 3323.6
$ws.Range("I68").Value = 3504.8572
$ws.Range("K68").Value = 3504.8572
$ws.Range("M68").Value = -2755.8572

$ws.Range("H71").Value = 3323.6
$ws.Range("I71").Value = 3504.8572
$ws.Range("K71").Value = 17524.286
$ws.Range("M71").Value = -13780.286

$ws.Range("H82").Value = 2523.625
$ws.Range("J82").Value = 2011.4286
$ws.Range("L82").Value = 2011.4286
$ws.Range("N82").Value = -2733.4286

$ws.Range("H85").Value = 2523.625
$ws.Range("J85").Value = 2011.4286
$ws.Range("L85").Value = 2011.4286
$ws.Range("N85").Value = -4507.4286

$ws.Range("H122").Value = 560695
$ws.Range("I122").Value = 773285.9
$ws.Range("K122").Value = 2319857.7
$ws.Range("M122").Value = -2317407.7

$ws.Range("H132").Value = 7276.4443
$ws.Range("J132").Value = 8297.799999999999
$ws.Range("L132").Value = 24893.4
$ws.Range("N132").Value = -29953.4

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 9573.666999999999
$ws.Range("I62").Value = 4248.3335
$ws.Range("K62").Value = 4248.3335
$ws.Range("M62").Value = -3624.3335

$ws.Range("H65").Value = 9573.666999999999
$ws.Range("I65").Value = 4248.3335
$ws.Range("K65").Value = 21241.6675
$ws.Range("M65").Value = -18121.6675

$ws.Range("H113").Value = 1341.5
$ws.Range("I113").Value = 827.3158
$ws.Range("J113").Value = 4598
$ws.Range("K113").Value = 2481.9474
$ws.Range("L113").Value = 13794
$ws.Range("M113").Value = -311.9474
$ws.Range("N113").Value = -18134

$ws.Range("H122").Value = 1975.0714
$ws.Range("I122").Value = 1712.5834
$ws.Range("K122").Value = 5137.7502
$ws.Range("M122").Value = -2687.7502

$ws.Range("H126").Value = 19000.428
$ws.Range("I126").Value = 28626
$ws.Range("K126").Value = 85878
$ws.Range("M126").Value = -83408

$ws.Range("H132").Value = 2739.25
$ws.Range("I132").Value = 2416.2856
$ws.Range("K132").Value = 7248.8568
$ws.Range("M132").Value = -4718.8568

$ws.Range("H136").Value = 2928.52
$ws.Range("I136").Value = 2956.3333
$ws.Range("K136").Value = 8868.999899999999
$ws.Range("M136").Value = -6318.999899999999
